$d = $word.ActiveDocument

$replacements = @(
    @("38×11=418", "20×51=1020"),
    @("62×16=992", "29×95=2755"),
    @("12×97=1164", "47×73=3431"),
    @("55×28=1540", "51×62=3162"),
    @("63×29=1827", "14×74=1036"),
    @("44×83=3652", "30×37=1110"),
    @("21×14=294", "79×69=5451"),
    @("57×96=5472", "44×17=748"),
    @("23×41=943", "56×43=2408"),
    @("99×22=2178", "57×14=798"),
    @("68×72=4896", "63×30=1890"),
    @("56×85=4760", "87×47=4089"),
    @("50×41=2050", "46×17=782"),
    @("76×41=3116", "16×63=1008"),
    @("40×82=3280", "31×34=1054"),
    @("58×67=3886", "77×97=7469"),
    @("99×94=9306", "68×23=1564"),
    @("11×69=759", "87×17=1479"),
    @("41×95=3895", "76×55=4180"),
    @("77×16=1232", "29×55=1595"),
    @("50×58=2900", "49×40=1960"),
    @("52×74=3848", "65×50=3250"),
    @("56×91=5096", "75×64=4800"),
    @("75×45=3375", "51×48=2448"),
    @("28×75=2100", "20×73=1460")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
